$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coureurs")

# Most recent race results entered in column G (Saudi-Arabië)
$ws.Range("G2").Value = 19
$ws.Range("G3").Value = 25
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 12
$ws.Range("G7").Value = 10
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 2
$ws.Range("G15").Value = 15
$ws.Range("G16").Value = 1

# Update the view/selection to match the sheet's latest state
$ws.Activate()
$ws.Range("G7").Select()
